$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 183, pushing existing rows 183-188 down to 184-189.
$ws.Rows.Item(183).Insert()

# Populate the newly inserted row 183 with the new weekly record.
$ws.Cells.Item(183, 1).Value = 4
$ws.Cells.Item(183, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(183, 3).Value = "Los Lagos"
$ws.Cells.Item(183, 4).Value = 44747
$ws.Cells.Item(183, 5).Value = 10
$ws.Cells.Item(183, 6).Value = 100112009
$ws.Cells.Item(183, 7).Value = "Acelga"
$ws.Cells.Item(183, 8).Value = "Sin especificar"
$ws.Cells.Item(183, 9).Value = "Primera"
$ws.Cells.Item(183, 10).Value = 90
$ws.Cells.Item(183, 11).Value = 12000
$ws.Cells.Item(183, 12).Value = 12000
$ws.Cells.Item(183, 13).Value = 12000
$ws.Cells.Item(183, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(183, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(183, 16).Value = 1000
$ws.Cells.Item(183, 17).Value = 12
$ws.Cells.Item(183, 18).Value = "Hortaliza"

Write-Host "row 183 inserted"
